# Daily attendance processing - 2025-10-20 07:20:54
# Normalize the "Recorded By" (column G) values so that entries recorded as
# "System, <email>" are rewritten as "<email>, System" (moving "System" to
# the end of the list). Entries that already have "System" at the end, or
# that have more/other combinations, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.StartsWith("System, ")) {
        $rest = $value.Substring(8)
        $cell.Value = "$rest, System"
    }
}
